$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header "vartio" stays a string, new numeric headers 449 / 450
$ws.Range("A1").Value = "vartio"
$ws.Range("B1").Value = 449
$ws.Range("C1").Value = 450

# Row 2 / Row 3: column C becomes numeric instead of text
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 2

# Update the remembered selection in the sheet view
$ws.Range("B8").Select()
